$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handoff
#
# The source file "22b44343-3293-45c8-bb58-a349678de8aa.md" was re-handed-off
# and now carries a new guid-named file + new package hash, with fresh
# handoff timestamps. The previous "7a8d8ed1-...md" entry (which had failed
# transform) is gone from the report, so every sheet loses that row and the
# ".localization-config" row moves up to take its place.
# ---------------------------------------------------------------------------

$newMd        = "744ee916-8fe6-421b-8f6f-1069caca09ba.md"
$newZhXlf     = "744ee916-8fe6-421b-8f6f-1069caca09ba.414792c94634192d6c7ec99ae69c61e42dc6dc8f.zh-cn.xlf"
$newDeXlf     = "744ee916-8fe6-421b-8f6f-1069caca09ba.414792c94634192d6c7ec99ae69c61e42dc6dc8f.de-de.xlf"
$newZhStamp   = "2016-02-06 04:16:52"
$newDeStamp   = "2016-02-06 04:17:02"
$cfgName      = ".localization-config"
$cfgStatus    = "Not to be localized"

# ===========================================================================
# Sheet 1: "Overview"
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

# Row 2: rename the handed-off source file.
$ws.Range("A2").Value2 = $newMd
$links = @($ws.Hyperlinks)
$links[0].TextToDisplay = $newMd

# Row 3 used to be the "7a8d8ed1-....md" / "Handoff transform failed" entry.
# That file no longer appears in the report, so row 3 is overwritten with
# what used to be row 4's content (".localization-config").
$ws.Range("A3").Value2 = $cfgName
$ws.Range("B3").Value2 = $cfgStatus
$ws.Range("C3").Value2 = $cfgStatus
$links = @($ws.Hyperlinks)
$links[1].TextToDisplay = $cfgName

# The hyperlink that used to anchor row 4 is now redundant.
$links = @($ws.Hyperlinks)
$links[2].Delete()

# Remove the now-duplicate trailing row.
$ws.Rows("4").Delete()

# ===========================================================================
# Sheet 2: "zh-cn"
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2: rename source file, target xlf package, and handoff timestamp.
$ws.Range("A2").Value2 = $newMd
$ws.Range("C2").Value2 = $newZhXlf
$ws.Range("D2").Value2 = $newZhStamp
$links = @($ws.Hyperlinks)
$links[0].TextToDisplay = $newMd
$links[1].TextToDisplay = $newZhXlf

# Row 3 ("7a8d8ed1-....md") is replaced by what used to be row 4.
$ws.Range("A3").Value2 = $cfgName
$ws.Range("B3").Value2 = $cfgStatus
$links = @($ws.Hyperlinks)
$links[2].TextToDisplay = $cfgName

# Drop the hyperlink that used to anchor row 4, then the row itself.
$links = @($ws.Hyperlinks)
$links[3].Delete()
$ws.Rows("4").Delete()

# ===========================================================================
# Sheet 3: "de-de"
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

# Row 2: rename source file, target xlf package, and handoff timestamp.
$ws.Range("A2").Value2 = $newMd
$ws.Range("C2").Value2 = $newDeXlf
$ws.Range("D2").Value2 = $newDeStamp
$links = @($ws.Hyperlinks)
$links[0].TextToDisplay = $newMd
$links[1].TextToDisplay = $newDeXlf

# Row 3 ("7a8d8ed1-....md") is replaced by what used to be row 4.
$ws.Range("A3").Value2 = $cfgName
$ws.Range("B3").Value2 = $cfgStatus
$links = @($ws.Hyperlinks)
$links[2].TextToDisplay = $cfgName

# Drop the hyperlink that used to anchor row 4, then the row itself.
$links = @($ws.Hyperlinks)
$links[3].Delete()
$ws.Rows("4").Delete()
